$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

# Insert two new rows after row 90 (CookieEnhanceStarReq's CostObj row),
# shifting everything below down by 2.
$ws.Rows("91:92").Insert()

# Row 90: replace CostObj/CostObjPacket with UsedSoulStone/int (count stays 5)
$ws.Range("C90").Value = "UsedSoulStone"
$ws.Range("D90").Value = "int"

# Row 91 (new): BefAccSoulStone/int, count 6
$ws.Range("B91").Value = "CookieEnhanceStarReq"
$ws.Range("C91").Value = "BefAccSoulStone"
$ws.Range("D91").Value = "int"
$ws.Range("E91").Value = 6

# Row 92 (new): AftAccSoulStone/int, count 7
$ws.Range("B92").Value = "CookieEnhanceStarReq"
$ws.Range("C92").Value = "AftAccSoulStone"
$ws.Range("D92").Value = "int"
$ws.Range("E92").Value = 7

# The old "ChgObj/ChgObjPacket" detail row that followed CookieEnhanceStarRes/Cookie
# has moved down to row 95 after the insert above; remove it entirely.
$ws.Rows("95:95").Delete()

# Update sheet view/selection to match the authored state.
$ws.Range("A92").Select()
$ws.Application.ActiveWindow.ScrollRow = 92
$ws.Range("C101").Select()

$wb.Windows.Item(1).WindowState = -4143
